$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the C column values for rows 1 through 71 (the "ResultValue" helper column
# that is no longer needed in the final results sheet).
$ws.Range("C1:C71").ClearContents()

# Update the view: reset scroll position (clears the saved topLeftCell="A45") and
# set a fresh single-cell selection at H22 (matches the final saved view state).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H22").Select()
